$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("magnesian/calc-alkalic/peraluminous", 893),
    @("magnesian/calc-alkalic/metaluminous", 630),
    @("ferroan/calc-alkalic/peraluminous", 605),
    @("magnesian/alkali-calcic/peraluminous", 492),
    @("ferroan/alkali-calcic/peraluminous", 397),
    @("magnesian/calcic/peraluminous", 346),
    @("magnesian/alkali-calcic/metaluminous", 330),
    @("magnesian/calcic/metaluminous", 179),
    @("ferroan/alkali-calcic/metaluminous", 177),
    @("ferroan/calcic/peraluminous", 119),
    @("ferroan/calc-alkalic/metaluminous", 104),
    @("magnesian/alkalic/metaluminous", 77),
    @("magnesian/alkalic/peraluminous", 67),
    @("ferroan/alkalic/peraluminous", 59),
    @("ferroan/alkalic/metaluminous", 44),
    @("ferroan/alkali-calcic/peralkaline", 42),
    @("magnesian/alkali-calcic/peralkaline", 31),
    @("ferroan/alkalic/peralkaline", 22),
    @("ferroan/calcic/metaluminous", 20),
    @("magnesian/alkalic/peralkaline", 15),
    @("ferroan/calc-alkalic/peralkaline", 5),
    @("magnesian/calc-alkalic/peralkaline", 4),
    @("magnesian/calcic/peralkaline", 1)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    # Copy style from existing row 2 (A column style s="1") to new rows beyond 13
    if ($row -gt 13) {
        $ws.Range("A2").Copy() | Out-Null
        $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    }
    $row++
}
